$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$win = $excel.ActiveWindow

$win.FreezePanes = $false
$ws.Range("A3").Select()
$win.FreezePanes = $true
Write-Host "done"
